$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.661.57"
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.533.23"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.34"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.59"
$ws.Range("E6").Value = "  -1.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  -0.96%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -2.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.45"
$ws.Range("E10").Value = "  -1.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0804"
$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.110"
$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.27"
$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.924.16"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.41"
$ws.Range("E15").Value = "  -3.87%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.535.70"
$ws.Range("E16").Value = "  -2.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.812"
$ws.Range("E17").Value = "  -3.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.660.38"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.64"
$ws.Range("E19").Value = "  -1.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.28"
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0948"
$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.68"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.19"
$ws.Range("E23").Value = "  -2.60%  "

$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.01"
$ws.Range("E25").Value = "  -2.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.51"
$ws.Range("E27").Value = "  -3.97%  "

$ws.Range("E28").Value = "  -2.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.09"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.15"
$ws.Range("E30").Value = "  -5.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.85"
$ws.Range("E31").Value = "  +2.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.76"
$ws.Range("E32").Value = "  -0.35%  "

$ws.Range("E33").Value = "  +1.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0790"
$ws.Range("E34").Value = "  -2.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.64"
$ws.Range("E35").Value = "  +1.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.15"
$ws.Range("E36").Value = "  -4.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.95"
$ws.Range("E37").Value = "  -6.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.59"
$ws.Range("E38").Value = "  -3.33%  "

$ws.Range("E39").Value = "  -1.16%  "

$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.12"
$ws.Range("E41").Value = "  -1.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.84"
$ws.Range("E42").Value = "  -3.34%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.32"
$ws.Range("E44").Value = "  +0.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0298"
$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.990.12"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.08"
$ws.Range("E47").Value = "  +0.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.780.36"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.190"
$ws.Range("E49").Value = "  -1.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.11"
$ws.Range("E50").Value = "  -2.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.89"
$ws.Range("E51").Value = "  -2.05%  "
